$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "287.51"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.72%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "29.20"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.06%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.089"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.16%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06686"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.01%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.342"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.53%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.407"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.20%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.373"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.57%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9195"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.81%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1585"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.88%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06805"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "4.61%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07569"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.34%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02929"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.73%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.08968"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.15%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001586"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.96%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04508"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.12%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0006439"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.22%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006295"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "4.32%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.26%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.07%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.42%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.062"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.22%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1582"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.84%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001193"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.89%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004111"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-5.05%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "1.66%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001617"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-1.10%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04262"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.89%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006721"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.83%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1238"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.37%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002210"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.21%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01341"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "12.92%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005679"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "5.50%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.81%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01306"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-29.42%"
